# Update plainLiterals -> xsd:string (rdf:langString when a language list is set)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("class-based constraints")

# Rows where column M currently holds "rdf:plainLiteral".
# If the row also specifies sh:languageIn (column T), the datatype should
# become rdf:langString; otherwise it becomes xsd:string.
$rows = 8..39
foreach ($r in $rows) {
    $mCell = $ws.Cells.Item($r, 13)   # column M
    if ($mCell.Value2 -eq "rdf:plainLiteral") {
        $tCell = $ws.Cells.Item($r, 20)   # column T (sh:languageIn)
        if ($tCell.Value2) {
            $mCell.Value2 = "rdf:langString"
        } else {
            $mCell.Value2 = "xsd:string"
        }
    }
}
